$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(45743,45743.01041666666,45743.02083333334,45743.03125,45743.04166666666,45743.05208333334,45743.0625,45743.07291666666,45743.08333333334,45743.09375,45743.10416666666,45743.11458333334,45743.125,45743.13541666666,45743.14583333334,45743.15625,45743.16666666666,45743.17708333334,45743.1875,45743.19791666666,45743.20833333334,45743.21875,45743.22916666666,45743.23958333334,45743.25,45743.26041666666,45743.27083333334,45743.28125,45743.29166666666,45743.30208333334,45743.3125,45743.32291666666,45743.33333333334,45743.34375,45743.35416666666,45743.36458333334,45743.375,45743.38541666666,45743.39583333334,45743.40625,45743.41666666666,45743.42708333334,45743.4375,45743.44791666666,45743.45833333334,45743.46875,45743.47916666666,45743.48958333334,45743.5,45743.51041666666,45743.52083333334,45743.53125,45743.54166666666,45743.55208333334,45743.5625,45743.57291666666,45743.58333333334,45743.59375,45743.60416666666,45743.61458333334,45743.625,45743.63541666666,45743.64583333334,45743.65625,45743.66666666666,45743.67708333334,45743.6875,45743.69791666666,45743.70833333334,45743.71875,45743.72916666666,45743.73958333334,45743.75,45743.76041666666,45743.77083333334,45743.78125,45743.79166666666,45743.80208333334,45743.8125,45743.82291666666,45743.83333333334,45743.84375,45743.85416666666,45743.86458333334,45743.875,45743.88541666666,45743.89583333334,45743.90625,45743.91666666666,45743.92708333334,45743.9375,45743.94791666666,45743.95833333334,45743.96875,45743.97916666666,45743.98958333334,45744,45744.01041666666,45744.02083333334,45744.03125,45744.04166666666,45744.05208333334,45744.0625,45744.07291666666,45744.08333333334,45744.09375,45744.10416666666,45744.11458333334,45744.125,45744.13541666666,45744.14583333334,45744.15625,45744.16666666666,45744.17708333334,45744.1875,45744.19791666666,45744.20833333334,45744.21875,45744.22916666666,45744.23958333334,45744.25,45744.26041666666,45744.27083333334,45744.28125,45744.29166666666,45744.30208333334,45744.3125,45744.32291666666,45744.33333333334,45744.34375,45744.35416666666,45744.36458333334,45744.375,45744.38541666666,45744.39583333334,45744.40625,45744.41666666666,45744.42708333334,45744.4375,45744.44791666666,45744.45833333334,45744.46875,45744.47916666666,45744.48958333334,45744.5,45744.51041666666,45744.52083333334,45744.53125,45744.54166666666,45744.55208333334,45744.5625,45744.57291666666,45744.58333333334,45744.59375,45744.60416666666,45744.61458333334,45744.625,45744.63541666666,45744.64583333334,45744.65625,45744.66666666666,45744.67708333334,45744.6875,45744.69791666666,45744.70833333334,45744.71875,45744.72916666666,45744.73958333334,45744.75,45744.76041666666,45744.77083333334,45744.78125,45744.79166666666,45744.80208333334,45744.8125,45744.82291666666,45744.83333333334,45744.84375,45744.85416666666,45744.86458333334,45744.875,45744.88541666666,45744.89583333334,45744.90625,45744.91666666666,45744.92708333334,45744.9375,45744.94791666666,45744.95833333334,45744.96875,45744.97916666666,45744.98958333334)
$bValues = @(150,148,148,147,148,150,150,149,142,147,147,147,147,148,146,145,145,145,146,153,209,219,218,217,283,292,294,307,413,428,435,431,372,366,360,358,247,249,250,247,232,230,220,230,150,146,146,158,149,152,212,210,219,220,218,215,181,191,181,192,262,195,182,173,187,188,187,190,332,339,345,353,580,596,598,606,575,474,463,459,349,339,427,422,205,229,249,247,137,109,107,100,56,47,41,41,38,38,38,38,38,38,38,38,38,38,37,38,41,41,41,41,41,41,42,45,86,112,119,127,235,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$startRow = 2
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $aValues[$i]
    $ws.Cells.Item($r, 2).Value = $bValues[$i]
}
